$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('ADBE', 'Adobe Inc.', 'Information Technology', 'Application Software'),
    @('ADP', 'ADP', 'Information Technology', 'Data Processing & Outsourced Services'),
    @('ABNB', 'Airbnb', 'Consumer Discretionary', 'Internet & Direct Marketing Retail'),
    @('ALGN', 'Align Technology', 'Health Care', 'Health Care Supplies'),
    @('GOOGL', 'Alphabet Inc. (Class A)', 'Communication Services', 'Interactive Media & Services'),
    @('GOOG', 'Alphabet Inc. (Class C)', 'Communication Services', 'Interactive Media & Services'),
    @('AMZN', 'Amazon', 'Consumer Discretionary', 'Internet & Direct Marketing Retail'),
    @('AMD', 'Advanced Micro Devices Inc.', 'Information Technology', 'Semiconductors'),
    @('AEP', 'American Electric Power', 'Utilities', 'Electric Utilities'),
    @('AMGN', 'Amgen', 'Health Care', 'Biotechnology'),
    @('ADI', 'Analog Devices', 'Information Technology', 'Semiconductors'),
    @('ANSS', 'Ansys', 'Information Technology', 'Application Software'),
    @('AAPL', 'Apple Inc.', 'Information Technology', 'Technology Hardware, Storage & Peripherals'),
    @('AMAT', 'Applied Materials', 'Information Technology', 'Semiconductor Equipment'),
    @('ASML', 'ASML Holding', 'Information Technology', 'Semiconductor Equipment'),
    @('AZN', 'AstraZeneca', 'Health Care', 'Pharmaceuticals'),
    @('TEAM', 'Atlassian', 'Information Technology', 'Application Software'),
    @('ADSK', 'Autodesk', 'Information Technology', 'Application Software'),
    @('BKR', 'Baker Hughes', 'Energy', 'Oil & Gas Equipment & Services'),
    @('BIIB', 'Biogen', 'Health Care', 'Biotechnology'),
    @('BKNG', 'Booking Holdings', 'Consumer Discretionary', 'Internet & Direct Marketing Retail'),
    @('AVGO', 'Broadcom Inc.', 'Information Technology', 'Semiconductors'),
    @('CDNS', 'Cadence Design Systems', 'Information Technology', 'Application Software'),
    @('CHTR', 'Charter Communications', 'Communication Services', 'Cable & Satellite'),
    @('CTAS', 'Cintas', 'Industrials', 'Diversified Support Services'),
    @('CSCO', 'Cisco', 'Information Technology', 'Communications Equipment'),
    @('CTSH', 'Cognizant', 'Information Technology', 'IT Consulting & Other Services'),
    @('CMCSA', 'Comcast', 'Communication Services', 'Cable & Satellite'),
    @('CEG', 'Constellation Energy', 'Utilities', 'Multi-Utilities'),
    @('CPRT', 'Copart', 'Industrials', 'Diversified Support Services'),
    @('CSGP', 'CoStar Group', 'Industrials', 'Research & Consulting Services'),
    @('COST', 'Costco', 'Consumer Staples', 'Hypermarkets & Super Centers'),
    @('CRWD', 'CrowdStrike', 'Information Technology', 'Application Software'),
    @('CSX', 'CSX Corporation', 'Industrials', 'Railroads'),
    @('DDOG', 'Datadog', 'Information Technology', 'Application Software'),
    @('DXCM', 'DexCom', 'Health Care', 'Health Care Equipment'),
    @('FANG', 'Diamondback Energy', 'Energy', 'Oil & Gas Exploration & Production'),
    @('DLTR', 'Dollar Tree', 'Consumer Discretionary', 'General Merchandise Stores'),
    @('EBAY', 'eBay', 'Consumer Discretionary', 'Internet & Direct Marketing Retail'),
    @('EA', 'Electronic Arts', 'Communication Services', 'Interactive Home Entertainment'),
    @('ENPH', 'Enphase Energy', 'Information Technology', 'Electronic Components'),
    @('EXC', 'Exelon', 'Utilities', 'Multi-Utilities'),
    @('FAST', 'Fastenal', 'Industrials', 'Building Products'),
    @('FTNT', 'Fortinet', 'Information Technology', 'Systems Software'),
    @('GEHC', 'GE HealthCare', 'Health Care', 'Health Care Technology'),
    @('GILD', 'Gilead Sciences', 'Health Care', 'Biotechnology'),
    @('GFS', 'GlobalFoundries', 'Information Technology', 'Semiconductors'),
    @('HON', 'Honeywell', 'Industrials', 'Industrial Conglomerates'),
    @('IDXX', 'Idexx Laboratories', 'Health Care', 'Health Care Equipment'),
    @('ILMN', 'Illumina, Inc.', 'Health Care', 'Life Sciences Tools & Services'),
    @('INTC', 'Intel', 'Information Technology', 'Semiconductors'),
    @('INTU', 'Intuit', 'Information Technology', 'Application Software'),
    @('ISRG', 'Intuitive Surgical', 'Health Care', 'Health Care Equipment'),
    @('JD', 'JD.com', 'Consumer Discretionary', 'Internet & Direct Marketing Retail'),
    @('KDP', 'Keurig Dr Pepper', 'Consumer Staples', 'Soft Drinks'),
    @('KLAC', 'KLA Corporation', 'Information Technology', 'Semiconductor Equipment'),
    @('KHC', 'Kraft Heinz', 'Consumer Staples', 'Packaged Foods & Meats'),
    @('LRCX', 'Lam Research', 'Information Technology', 'Semiconductor Equipment'),
    @('LCID', 'Lucid Motors', 'Consumer Discretionary', 'Automobile Manufacturers'),
    @('LULU', 'Lululemon', 'Consumer Discretionary', 'Apparel, Accessories & Luxury Goods'),
    @('MAR', 'Marriott International', 'Consumer Discretionary', 'Hotels, Resorts & Cruise Lines'),
    @('MRVL', 'Marvell Technology', 'Information Technology', 'Application Software'),
    @('MELI', 'MercadoLibre', 'Consumer Discretionary', 'Internet & Direct Marketing Retail'),
    @('META', 'Meta Platforms', 'Communication Services', 'Interactive Media & Services'),
    @('MCHP', 'Microchip Technology', 'Information Technology', 'Semiconductors'),
    @('MU', 'Micron Technology', 'Information Technology', 'Semiconductors'),
    @('MSFT', 'Microsoft', 'Information Technology', 'Systems Software'),
    @('MRNA', 'Moderna', 'Health Care', 'Biotechnology'),
    @('MDLZ', 'Mondelēz International', 'Consumer Staples', 'Packaged Foods & Meats'),
    @('MNST', 'Monster Beverage', 'Consumer Staples', 'Soft Drinks'),
    @('NFLX', 'Netflix', 'Communication Services', 'Movies & Entertainment'),
    @('NVDA', 'Nvidia', 'Information Technology', 'Semiconductors'),
    @('NXPI', 'NXP', 'Information Technology', 'Semiconductors'),
    @('ORLY', 'O''Reilly Automotive', 'Consumer Discretionary', 'Specialty Stores'),
    @('ODFL', 'Old Dominion Freight Line', 'Industrials', 'Trucking'),
    @('ON', 'Onsemi', 'Information Technology', 'Semiconductors'),
    @('PCAR', 'Paccar', 'Industrials', 'Construction Machinery & Heavy Trucks'),
    @('PANW', 'Palo Alto Networks', 'Information Technology', 'Application Software'),
    @('PAYX', 'Paychex', 'Information Technology', 'Data Processing & Outsourced Services'),
    @('PYPL', 'PayPal', 'Information Technology', 'Data Processing & Outsourced Services'),
    @('PDD', 'PDD Holdings', 'Consumer Discretionary', 'Internet & Direct Marketing Retail'),
    @('PEP', 'PepsiCo', 'Consumer Staples', 'Soft Drinks'),
    @('QCOM', 'Qualcomm', 'Information Technology', 'Semiconductors'),
    @('REGN', 'Regeneron', 'Health Care', 'Biotechnology'),
    @('ROST', 'Ross Stores', 'Consumer Discretionary', 'Apparel Retail'),
    @('SGEN', 'Seagen', 'Health Care', 'Biotechnology'),
    @('SIRI', 'Sirius XM', 'Communication Services', 'Broadcasting'),
    @('SBUX', 'Starbucks', 'Consumer Discretionary', 'Restaurants'),
    @('SNPS', 'Synopsys', 'Information Technology', 'Application Software'),
    @('TMUS', 'T-Mobile US', 'Communication Services', 'Wireless Telecommunication Services'),
    @('TSLA', 'Tesla, Inc.', 'Consumer Discretionary', 'Automobile Manufacturers'),
    @('TXN', 'Texas Instruments', 'Information Technology', 'Semiconductors'),
    @('TTD', 'The Trade Desk', 'Communication Services', '')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($row[3] -eq "") {
        $ws.Cells.Item($r, 4).Value = $null
    } else {
        $ws.Cells.Item($r, 4).Value = $row[3]
    }
}
